$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text: volume number and report week dates (rich-text runs) ---
# "Volume 31   Number  10" -> "...  11"  (characters 21-22 hold "10")
$ws.Range("A8").Characters(21, 2).Text = "11"

# "Report Covering the Week  3/4/2024  Through  3/10/2024"
# Edit the right-hand date first so the left-hand edit's length change
# (8 chars -> 9 chars) doesn't shift the second run's start offset.
$ws.Range("C9").Characters(46, 9).Text = "3/17/2024"
$ws.Range("C9").Characters(27, 8).Text = "3/11/2024"

# --- Weekly crime-stat table body (rows 15-30) ---
# Style 15 = number "#,##0", Style 16 = percent-like "#,##0.0", Style 14 = text
# placeholder (shared strings "0" / "***.*") used when data is N/A. When a
# cell's value type flips we borrow the destination format (via a
# formats-only paste from a stable donor cell of the same column class)
# instead of letting the paste clobber the numbers we just wrote.
$ws.Range("I14").Copy()
$ws.Range("C15").PasteSpecial(-4122)
$ws.Range("C15").Value = 1
$ws.Range("I14").Copy()
$ws.Range("D15").PasteSpecial(-4122)
$ws.Range("D15").Value = 1
$ws.Range("K14").Copy()
$ws.Range("E15").PasteSpecial(-4122)
$ws.Range("E15").Value = 0
$ws.Range("F15").Value = 2
$ws.Range("I14").Copy()
$ws.Range("G15").PasteSpecial(-4122)
$ws.Range("G15").Value = 1
$ws.Range("K14").Copy()
$ws.Range("H15").PasteSpecial(-4122)
$ws.Range("H15").Value = 100
$ws.Range("I15").Value = 2
$ws.Range("I14").Copy()
$ws.Range("J15").PasteSpecial(-4122)
$ws.Range("J15").Value = 1
$ws.Range("K14").Copy()
$ws.Range("K15").PasteSpecial(-4122)
$ws.Range("K15").Value = 100
$ws.Range("L15").Value = -33.333333333333
$ws.Range("M15").Value = 100
$ws.Range("N15").Value = -66.666666666666
$ws.Range("C16").Value = 2
$ws.Range("D16").Value = "'0"
$ws.Range("M14").Copy()
$ws.Range("D16").PasteSpecial(-4122)
$ws.Range("E16").Value = "***.*"
$ws.Range("M14").Copy()
$ws.Range("E16").PasteSpecial(-4122)
$ws.Range("F16").Value = 6
$ws.Range("G16").Value = 3
$ws.Range("H16").Value = 100
$ws.Range("I16").Value = 22
$ws.Range("K16").Value = -4.347826086956
$ws.Range("L16").Value = -18.518518518518
$ws.Range("M16").Value = -43.589743589743
$ws.Range("N16").Value = -80.869565217391
$ws.Range("C17").Value = 1
$ws.Range("D17").Value = 4
$ws.Range("F17").Value = 12
$ws.Range("H17").Value = -40
$ws.Range("I17").Value = 25
$ws.Range("J17").Value = 49
$ws.Range("K17").Value = -48.979591836734
$ws.Range("L17").Value = -58.333333333333
$ws.Range("M17").Value = -19.354838709677
$ws.Range("N17").Value = -21.875
$ws.Range("F18").Value = 1
$ws.Range("G18").Value = 6
$ws.Range("H18").Value = -83.333333333333
$ws.Range("L18").Value = -46.666666666666
$ws.Range("M18").Value = -82.978723404255
$ws.Range("N18").Value = -94.366197183098
$ws.Range("C19").Value = 8
$ws.Range("D19").Value = 4
$ws.Range("E19").Value = 100
$ws.Range("F19").Value = 26
$ws.Range("G19").Value = 18
$ws.Range("H19").Value = 44.444444444444
$ws.Range("I19").Value = 65
$ws.Range("J19").Value = 48
$ws.Range("K19").Value = 35.416666666666
$ws.Range("L19").Value = 18.181818181818
$ws.Range("M19").Value = 35.416666666666
$ws.Range("N19").Value = -16.666666666666
$ws.Range("D20").Value = 3
$ws.Range("E20").Value = -33.333333333333
$ws.Range("F20").Value = 6
$ws.Range("H20").Value = -60
$ws.Range("I20").Value = 28
$ws.Range("J20").Value = 31
$ws.Range("K20").Value = -9.677419354838
$ws.Range("L20").Value = -22.222222222222
$ws.Range("M20").Value = 33.333333333333
$ws.Range("N20").Value = -93.926247288503
$ws.Range("C21").Value = 14
$ws.Range("D21").Value = 12
$ws.Range("E21").Value = 16.666666666666
$ws.Range("F21").Value = 53
$ws.Range("G21").Value = 63
$ws.Range("H21").Value = -15.873015873015
$ws.Range("I21").Value = 151
$ws.Range("J21").Value = 166
$ws.Range("K21").Value = -9.036144578313
$ws.Range("L21").Value = -23.737373737373
$ws.Range("M21").Value = -19.251336898395
$ws.Range("N21").Value = -81.89448441247
$ws.Range("I14").Copy()
$ws.Range("D22").PasteSpecial(-4122)
$ws.Range("D22").Value = 1
$ws.Range("K14").Copy()
$ws.Range("E22").PasteSpecial(-4122)
$ws.Range("E22").Value = -100
$ws.Range("I14").Copy()
$ws.Range("G22").PasteSpecial(-4122)
$ws.Range("G22").Value = 1
$ws.Range("K14").Copy()
$ws.Range("H22").PasteSpecial(-4122)
$ws.Range("H22").Value = -100
$ws.Range("J22").Value = 2
$ws.Range("K22").Value = 100
$ws.Range("C23").Value = 2
$ws.Range("D23").Value = "'0"
$ws.Range("M14").Copy()
$ws.Range("D23").PasteSpecial(-4122)
$ws.Range("E23").Value = "***.*"
$ws.Range("M14").Copy()
$ws.Range("E23").PasteSpecial(-4122)
$ws.Range("F23").Value = 4
$ws.Range("G23").Value = 9
$ws.Range("H23").Value = -55.555555555555
$ws.Range("I23").Value = 19
$ws.Range("K23").Value = -9.523809523809
$ws.Range("L23").Value = 18.75
$ws.Range("M23").Value = 216.666666666667
$ws.Range("C24").Value = 7
$ws.Range("D24").Value = 12
$ws.Range("E24").Value = -41.666666666666
$ws.Range("F24").Value = 32
$ws.Range("G24").Value = 48
$ws.Range("H24").Value = -33.333333333333
$ws.Range("I24").Value = 105
$ws.Range("J24").Value = 142
$ws.Range("K24").Value = -26.056338028169
$ws.Range("L24").Value = -5.405405405405
$ws.Range("M24").Value = 19.318181818181
$ws.Range("C25").Value = 1
$ws.Range("D25").Value = 5
$ws.Range("E25").Value = -80
$ws.Range("F25").Value = 6
$ws.Range("H25").Value = -45.454545454545
$ws.Range("I25").Value = 27
$ws.Range("J25").Value = 32
$ws.Range("K25").Value = -15.625
$ws.Range("L25").Value = -46
$ws.Range("C26").Value = 11
$ws.Range("D26").Value = 13
$ws.Range("E26").Value = -15.384615384615
$ws.Range("F26").Value = 24
$ws.Range("G26").Value = 26
$ws.Range("H26").Value = -7.692307692307
$ws.Range("I26").Value = 70
$ws.Range("J26").Value = 70
$ws.Range("K26").Value = 0
$ws.Range("L26").Value = 11.111111111111
$ws.Range("M26").Value = -16.666666666666
$ws.Range("I14").Copy()
$ws.Range("C27").PasteSpecial(-4122)
$ws.Range("C27").Value = 1
$ws.Range("I14").Copy()
$ws.Range("D27").PasteSpecial(-4122)
$ws.Range("D27").Value = 1
$ws.Range("K14").Copy()
$ws.Range("E27").PasteSpecial(-4122)
$ws.Range("E27").Value = 0
$ws.Range("F27").Value = 2
$ws.Range("I14").Copy()
$ws.Range("G27").PasteSpecial(-4122)
$ws.Range("G27").Value = 1
$ws.Range("K14").Copy()
$ws.Range("H27").PasteSpecial(-4122)
$ws.Range("H27").Value = 100
$ws.Range("I27").Value = 2
$ws.Range("I14").Copy()
$ws.Range("J27").PasteSpecial(-4122)
$ws.Range("J27").Value = 1
$ws.Range("K14").Copy()
$ws.Range("K27").PasteSpecial(-4122)
$ws.Range("K27").Value = 100
$ws.Range("L27").Value = -71.428571428571
$ws.Range("I14").Copy()
$ws.Range("C28").PasteSpecial(-4122)
$ws.Range("C28").Value = 1
$ws.Range("E28").Value = 0
$ws.Range("G28").Value = 3
$ws.Range("H28").Value = -33.333333333333
$ws.Range("I28").Value = 8
$ws.Range("J28").Value = 4
$ws.Range("K28").Value = 100
$ws.Range("L28").Value = 14.285714285714
$ws.Range("C29").Value = "'0"
$ws.Range("M14").Copy()
$ws.Range("C29").PasteSpecial(-4122)
$ws.Range("C30").Value = "'0"
$ws.Range("M14").Copy()
$ws.Range("C30").PasteSpecial(-4122)

$ws.Application.CutCopyMode = $false
